$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 376
$ws1.Range("F3").Value = 811
$ws1.Range("F4").Value = 280
$ws1.Range("F5").Value = 956
$ws1.Range("F6").Value = 2297

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 376
$ws4.Range("F3").Value = 811
$ws4.Range("F4").Value = 280
$ws4.Range("F7").Value = 956
$ws4.Range("F8").Value = 2297
